# Update Bloc 1 - Euler
# Remove the "ZoneTexte 5" textbox (the astrofrog.github.io URL caption)
# from slide 4 of the presentation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "ZoneTexte 5") {
        $sh.Delete()
    }
}
